$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (prices + 1h volume %) for cryptos.xlsx.
# Source values in the sheet are plain text (inline strings), including
# numeric-looking prices and "%"-suffixed change figures, so each cell is
# first switched to Text number format before the new literal is written;
# this prevents Excel from auto-converting them into numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.98%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.11%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.996"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.16%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07817"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.06%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.188"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-8.04%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.035"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.04%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.86%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9145"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.24%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09738"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.71%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1886"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.69%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08649"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.18%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03537"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.70%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09964"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.67%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001493"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.64%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005717"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.29%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.25%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.071"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.60%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.21%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.04%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.760"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "10.54%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.47%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.99%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004794"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.00%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.34%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01761"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.87%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04741"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.35%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008050"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.69%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.13%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007670"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "9.62%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002161"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.31%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "9.86%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006073"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.14%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.14%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.14%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.14%"
